# "Generate Report for Handback" — localization-status.xlsx
#
# 1. The "Ready for handoff" status text (shared by every status cell in the
#    workbook) becomes "Handed back: in sync with en-US" now that the
#    handback has completed.
# 2. The zh-cn sheet's "Latest Handback DateTime" (H2/H3) is refreshed to the
#    actual handback timestamp.
# 3. The de-de sheet's "Latest Handback DateTime" (H2/H3) is refreshed too
#    (it previously held the zero-date placeholder).
# 4. The de-de sheet gains "Latest Target File" (F) / "Latest Handback File"
#    (G) columns populated + hyperlinked, mirroring the zh-cn sheet which
#    already reports its handback file.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

# --- 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$newStatus = "Handed back: in sync with en-US"

$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- 2. zh-cn handback datetime refreshed ---
$zhcn.Range("H2").Value = "2016-03-12 04:28:51"
$zhcn.Range("H3").Value = "2016-03-12 04:28:51"

# --- 3. de-de handback datetime refreshed (was the 0001-01-01 placeholder) ---
$dede.Range("H2").Value = "2016-03-12 04:28:56"
$dede.Range("H3").Value = "2016-03-12 04:28:56"

# --- 4. de-de: add Latest Target File (F) / Latest Handback File (G) ---
$targetFile   = "a.md"
$handbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("F2").Value = $targetFile
$dede.Range("G2").Value = $handbackFile
$dede.Range("F3").Value = $targetFile
$dede.Range("G3").Value = $handbackFile

$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/6a0226c4e6f30af78dbcbc8df767f9bef201ab00/e2e/a.md", "", "", $targetFile)
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7da374a8c94e9996e346c4b012e2f94a70604268/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", $handbackFile)
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ef0eec17eafa2c198d8b321fd853b40d13e5dda5/e2e/a.md", "", "", $targetFile)
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7da374a8c94e9996e346c4b012e2f94a70604268/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", $handbackFile)

# Match the hyperlink look-and-feel already used across the workbook.
$dede.Range("F2:G3").Style = "HyperLink"

Write-Host "Handback report generated: status text, handback timestamps and de-de target/handback columns updated."
